$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 133.33333
$ws.Range("I6").Value = 199
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 597
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = -485
$ws.Range("N6").Value = -230

$ws.Range("H18").Value = 174.26086
$ws.Range("I18").Value = 150.36363
$ws.Range("K18").Value = 150.36363
$ws.Range("M18").Value = 133.63637

$ws.Range("H55").Value = 341
$ws.Range("I55").Value = 1038
$ws.Range("J55").Value = 92.07143000000001
$ws.Range("K55").Value = 1038
$ws.Range("L55").Value = 92.07143000000001
$ws.Range("M55").Value = -824
$ws.Range("N55").Value = -520.07143

$ws.Range("H62").Value = 2850
$ws.Range("I62").Value = 2358.3333
$ws.Range("K62").Value = 2358.3333
$ws.Range("M62").Value = -1734.3333

$ws.Range("H65").Value = 2850
$ws.Range("I65").Value = 2358.3333
$ws.Range("K65").Value = 11791.6665
$ws.Range("M65").Value = -8671.666499999999

$ws.Range("H86").Value = 19749
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 19749
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws.Range("H137").Value = 69992
$ws.Range("I137").Value = 4522.375
$ws.Range("J137").Value = 144814.42
$ws.Range("K137").Value = 13567.125
$ws.Range("L137").Value = 434443.26
$ws.Range("M137").Value = -11017.125
$ws.Range("N137").Value = -439543.26

$ws.Range("H138").Value = 2302.92
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2302.92
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 6908.76
$ws.Range("N138").Value = -17188.76
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1189.0625
$ws.Range("I2").Value = 1175.619
$ws.Range("J2").Value = 1214.7273
$ws.Range("K2").Value = 1175.619
$ws.Range("L2").Value = 1214.7273
$ws.Range("M2").Value = -1062.619
$ws.Range("N2").Value = -1440.7273

$ws.Range("H45").Value = 2017
$ws.Range("I45").Value = 2075.2
$ws.Range("K45").Value = 2075.2
$ws.Range("M45").Value = -1698.2

$ws.Range("H110").Value = 922
$ws.Range("I110").Value = 757
$ws.Range("J110").Value = 1499.5
$ws.Range("K110").Value = 757
$ws.Range("L110").Value = 1499.5
$ws.Range("M110").Value = 1288
$ws.Range("N110").Value = -5589.5

$ws.Range("H116").Value = 1189.0625
$ws.Range("I116").Value = 1175.619
$ws.Range("J116").Value = 1214.7273
$ws.Range("K116").Value = 1175.619
$ws.Range("L116").Value = 1214.7273
$ws.Range("M116").Value = 1118.381
$ws.Range("N116").Value = -5802.7273

$ws.Range("H132").Value = 11797.184
$ws.Range("I132").Value = 1375.9
$ws.Range("J132").Value = 58114
$ws.Range("K132").Value = 4127.700000000001
$ws.Range("L132").Value = 174342
$ws.Range("M132").Value = -1597.700000000001
$ws.Range("N132").Value = -179402

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1189.0625
$ws.Range("I3").Value = 1175.619
$ws.Range("J3").Value = 1214.7273
$ws.Range("K3").Value = 1175.619
$ws.Range("L3").Value = 1214.7273
$ws.Range("M3").Value = -1061.619
$ws.Range("N3").Value = -1442.7273

$ws.Range("H86").Value = 1839.9642
$ws.Range("I86").Value = 1490.25
$ws.Range("J86").Value = 2714.25
$ws.Range("K86").Value = 1490.25
$ws.Range("L86").Value = 2714.25
$ws.Range("M86").Value = -367.25
$ws.Range("N86").Value = -4960.25

$ws.Range("H89").Value = 1839.9642
$ws.Range("I89").Value = 1490.25
$ws.Range("J89").Value = 2714.25
$ws.Range("K89").Value = 7451.25
$ws.Range("L89").Value = 13571.25
$ws.Range("M89").Value = -1835.25
$ws.Range("N89").Value = -24803.25

$ws.Range("H105").Value = 3133.7727
$ws.Range("I105").Value = 3053.7856
$ws.Range("J105").Value = 3273.75
$ws.Range("K105").Value = 3053.7856
$ws.Range("L105").Value = 3273.75
$ws.Range("M105").Value = -1306.7856
$ws.Range("N105").Value = -6767.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1615
$ws.Range("I16").Value = 1615
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1615
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1328
$ws.Range("N16").ClearContents()

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H41").Value = 5059
$ws.Range("I41").Value = 5059
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5059
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4631
$ws.Range("N41").ClearContents()

$ws.Range("H58").Value = 13162.024
$ws.Range("I58").Value = 1005.4286
$ws.Range("J58").Value = 39345.46
$ws.Range("K58").Value = 1005.4286
$ws.Range("L58").Value = 39345.46
$ws.Range("M58").Value = -802.4286
$ws.Range("N58").Value = -39751.46

$ws.Range("H99").Value = 16671008
$ws.Range("I99").Value = 3669.5
$ws.Range("J99").Value = 41672016
$ws.Range("K99").Value = 3669.5
$ws.Range("L99").Value = 41672016
$ws.Range("M99").Value = -2171.5
$ws.Range("N99").Value = -41675012

$ws.Range("H113").Value = 1615
$ws.Range("I113").Value = 1615
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1615
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 555
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 16671008
$ws.Range("I126").Value = 3669.5
$ws.Range("J126").Value = 41672016
$ws.Range("K126").Value = 11008.5
$ws.Range("L126").Value = 125016048
$ws.Range("M126").Value = -8538.5
$ws.Range("N126").Value = -125020988

$ws.Range("H132").Value = 14981.23
$ws.Range("I132").Value = 19840.852
$ws.Range("K132").Value = 59522.556
$ws.Range("M132").Value = -56992.556

$ws.Range("H134").Value = 1056.68
$ws.Range("I134").Value = 887.8261
$ws.Range("J134").Value = 1200.5186
$ws.Range("K134").Value = 2663.4783
$ws.Range("L134").Value = 3601.5558
$ws.Range("M134").Value = -128.4782999999998
$ws.Range("N134").Value = -8671.5558

$ws.Range("H136").Value = 13162.024
$ws.Range("I136").Value = 1005.4286
$ws.Range("J136").Value = 39345.46
$ws.Range("K136").Value = 3016.2858
$ws.Range("L136").Value = 118036.38
$ws.Range("M136").Value = -466.2857999999997
$ws.Range("N136").Value = -123136.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 200
$ws.Range("K7").Value = 600
$ws.Range("M7").Value = -488

$ws.Range("H14").Value = 262.6154
$ws.Range("I14").Value = 262.6154
$ws.Range("K14").Value = 787.8462000000001
$ws.Range("M14").Value = -614.8462000000001

$ws.Range("H80").Value = 36667
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 52500.5
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 157501.5
$ws.Range("M80").Value = -14064
$ws.Range("N80").Value = -159373.5

$ws.Range("H83").Value = 36667
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 52500.5
$ws.Range("K83").Value = 45000
$ws.Range("L83").Value = 472504.5
$ws.Range("M83").Value = -40320
$ws.Range("N83").Value = -481864.5

$ws.Range("H92").Value = 624.125
$ws.Range("I92").Value = 642.8570999999999
$ws.Range("J92").Value = 493
$ws.Range("K92").Value = 1928.5713
$ws.Range("L92").Value = 1479
$ws.Range("M92").Value = -680.5712999999998
$ws.Range("N92").Value = -3975

$ws.Range("H131").Value = 719.92
$ws.Range("J131").Value = 749.3736
$ws.Range("L131").Value = 2248.1208
$ws.Range("N131").Value = -12328.1208

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17742.715
$ws.Range("I70").Value = 27449.75
$ws.Range("J70").Value = 4800
$ws.Range("K70").Value = 27449.75
$ws.Range("L70").Value = 4800
$ws.Range("M70").Value = -27179.75
$ws.Range("N70").Value = -5340

$ws.Range("H73").Value = 17742.715
$ws.Range("I73").Value = 27449.75
$ws.Range("J73").Value = 4800
$ws.Range("K73").Value = 27449.75
$ws.Range("L73").Value = 4800
$ws.Range("M73").Value = -26513.75
$ws.Range("N73").Value = -6672

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 29999
$ws.Range("J36").Value = 29999
$ws.Range("L36").Value = 29999
$ws.Range("N36").Value = -31123

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3247669.5
$ws.Range("I107").Value = 1178.25
$ws.Range("J107").Value = 7576324.5
$ws.Range("K107").Value = 3534.75
$ws.Range("L107").Value = 22728973.5
$ws.Range("M107").Value = -1614.75
$ws.Range("N107").Value = -22732813.5

$ws.Range("H136").Value = 27028936
$ws.Range("I136").Value = 43480030
$ws.Range("J136").Value = 2136.3572
$ws.Range("K136").Value = 130440090
$ws.Range("L136").Value = 6409.071599999999
$ws.Range("M136").Value = -130437540
$ws.Range("N136").Value = -11509.0716
